$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nino34")

# --- Update existing row 32 values (D32:U32 revised forecast) ---
$ws.Range("D32").Value = -0.159
$ws.Range("E32").Value = -0.246
$ws.Range("F32").Value = -0.364
$ws.Range("G32").Value = -0.467
$ws.Range("H32").Value = -0.507
$ws.Range("I32").Value = -0.472
$ws.Range("J32").Value = -0.395
$ws.Range("K32").Value = -0.331
$ws.Range("L32").Value = -0.311
$ws.Range("M32").Value = -0.325
$ws.Range("N32").Value = -0.346
$ws.Range("O32").Value = -0.354
$ws.Range("P32").Value = -0.353
$ws.Range("Q32").Value = -0.372
$ws.Range("R32").Value = -0.421
$ws.Range("S32").Value = -0.471
$ws.Range("T32").Value = -0.484
$ws.Range("U32").Value = -0.454

# --- Update existing row 33 values (B33:U33 revised forecast) ---
$ws.Range("B33").Value = -0.306
$ws.Range("C33").Value = -0.33
$ws.Range("D33").Value = -0.417
$ws.Range("E33").Value = -0.548
$ws.Range("F33").Value = -0.652
$ws.Range("G33").Value = -0.669
$ws.Range("H33").Value = -0.605
$ws.Range("I33").Value = -0.505
$ws.Range("J33").Value = -0.412
$ws.Range("K33").Value = -0.346
$ws.Range("L33").Value = -0.301
$ws.Range("M33").Value = -0.275
$ws.Range("N33").Value = -0.258
$ws.Range("O33").Value = -0.247
$ws.Range("P33").Value = -0.249
$ws.Range("Q33").Value = -0.269
$ws.Range("R33").Value = -0.288
$ws.Range("S33").Value = -0.291
$ws.Range("T33").Value = -0.281
$ws.Range("U33").Value = -0.254

# --- Add new row 34 (2025-09 init month forecast) ---
$ws.Range("A34").Value = "2025-09"
$ws.Range("B34").Value = -0.42
$ws.Range("C34").Value = -0.478
$ws.Range("D34").Value = -0.599
$ws.Range("E34").Value = -0.7
$ws.Range("F34").Value = -0.725
$ws.Range("G34").Value = -0.666
$ws.Range("H34").Value = -0.56
$ws.Range("I34").Value = -0.451
$ws.Range("J34").Value = -0.366
$ws.Range("K34").Value = -0.305
$ws.Range("L34").Value = -0.269
$ws.Range("M34").Value = -0.25
$ws.Range("N34").Value = -0.238
$ws.Range("O34").Value = -0.237
$ws.Range("P34").Value = -0.25
$ws.Range("Q34").Value = -0.262
$ws.Range("R34").Value = -0.264
$ws.Range("S34").Value = -0.258
$ws.Range("T34").Value = -0.236
$ws.Range("U34").Value = -0.186

# Copy formatting from row 33 (same layout: A = bold/bordered label, B:U = 0.00 numeric)
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial(-4122)

$ws.Range("B33:U33").Copy()
$ws.Range("B34:U34").PasteSpecial(-4122)

$excel.CutCopyMode = $false
